$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values that changed
$priceUpdates = @{
    2 = "256.01"
    3 = "23.35"
    4 = "6.127"
    5 = "0.06059"
    6 = "6.676"
    7 = "3.450"
    8 = "1.322"
    9 = "0.7916"
    10 = "0.01277"
    11 = "0.1523"
    12 = "0.08088"
    13 = "0.03373"
    14 = "0.03111"
    15 = "0.09241"
    16 = "3.758"
    17 = "0.001680"
    18 = "0.04806"
    19 = "0.006262"
    20 = "0.005922"
    21 = "0.001099"
    22 = "0.0001512"
    23 = "3.676"
    24 = "2.206"
    25 = "0.3355"
    26 = "0.1271"
    27 = "0.0006499"
    40 = "0.04545"
    41 = "0.007057"
    42 = "0.1077"
    43 = "0.003372"
    44 = "0.009837"
    45 = "0.002999"
    46 = "0.00005928"
    47 = "0.00000000756"
    48 = "0.7058"
    49 = "0.09315"
    50 = "0.00002117"
    51 = "0.01018"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Update Hora (G) column from 8 to 9 for all data rows (2-51)
$horaRange = $ws.Range("G2:G51")
$horaRange.NumberFormat = "@"
$horaRange.Value = "9"
$horaRange.Style = "Normal"
